$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (36 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K15").Value = 235.38
$ws.Range("I15").Value = 78.45999999999999
$ws.Range("H15").Value = 78.45999999999999
$ws.Range("M15").Value = -66.38
$ws.Range("H33").Value = 296.46667
$ws.Range("K33").Value = 285.29166
$ws.Range("N33").Value = -799.16666
$ws.Range("J33").Value = 341.16666
$ws.Range("I33").Value = 285.29166
$ws.Range("M33").Value = -56.29165999999998
$ws.Range("L33").Value = 341.16666
$ws.Range("J53").Value = 667.8570999999999
$ws.Range("H53").Value = 550.4545000000001
$ws.Range("K53").Value = 345
$ws.Range("M53").Value = 292
$ws.Range("N53").Value = -1941.8571
$ws.Range("L53").Value = 667.8570999999999
$ws.Range("I53").Value = 345
$ws.Range("J125").Value = 757.8182
$ws.Range("L125").Value = 6820.3638
$ws.Range("H125").Value = 712.9231
$ws.Range("I125").Value = 466
$ws.Range("K125").Value = 4194
$ws.Range("M125").Value = -1734
$ws.Range("N125").Value = -11740.3638
$ws.Range("H129").Value = 979.7659
$ws.Range("N129").Value = -13112.557
$ws.Range("J129").Value = 1037.519
$ws.Range("L129").Value = 3112.557
$ws.Range("I138").Value = 2815
$ws.Range("K138").Value = 8445
$ws.Range("M138").Value = -3305
$ws.Range("N138").Value = -13354803.5
$ws.Range("J138").Value = 4448174.5
$ws.Range("L138").Value = 13344523.5
$ws.Range("H138").Value = 3128458.5

# --- Sheet: ARM (46 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("J5").Value = 600
$ws.Range("H5").Value = 493.33334
$ws.Range("L5").Value = 600
$ws.Range("I5").Value = 440
$ws.Range("M5").Value = -328
$ws.Range("K5").Value = 440
$ws.Range("N5").Value = -824
$ws.Range("M32").Value = -37385.26
$ws.Range("K32").Value = 37672.26
$ws.Range("H32").Value = 4313370
$ws.Range("L32").Value = 10459685
$ws.Range("J32").Value = 10459685
$ws.Range("I32").Value = 37672.26
$ws.Range("N32").Value = -10460259
$ws.Range("N61").Value = -101939.6
$ws.Range("J61").Value = 101515.6
$ws.Range("H61").Value = 28630242
$ws.Range("L61").Value = 101515.6
$ws.Range("K61").Value = 40041732
$ws.Range("I61").Value = 40041732
$ws.Range("M61").Value = -40041520
$ws.Range("L97").Value = 2545
$ws.Range("H97").Value = 2315847.8
$ws.Range("J97").Value = 2545
$ws.Range("K97").Value = 2718161.2
$ws.Range("N97").Value = -3537
$ws.Range("I97").Value = 2718161.2
$ws.Range("M97").Value = -2717665.2
$ws.Range("H122").Value = 3705855.5
$ws.Range("M122").Value = -2322.4375
$ws.Range("N122").Value = -23822803
$ws.Range("J122").Value = 7939301
$ws.Range("K122").Value = 4772.4375
$ws.Range("I122").Value = 1590.8125
$ws.Range("L122").Value = 23817903
$ws.Range("H132").Value = 16199039
$ws.Range("M132").Value = -75155192
$ws.Range("I132").Value = 25052574
$ws.Range("K132").Value = 75157722
$ws.Range("I136").Value = 40041732
$ws.Range("K136").Value = 120125196
$ws.Range("H136").Value = 28630242
$ws.Range("M136").Value = -120122646
$ws.Range("L136").Value = 304546.8
$ws.Range("N136").Value = -309646.8
$ws.Range("J136").Value = 101515.6

# --- Sheet: BSM (11 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I4").Value = 440
$ws.Range("L4").Value = 600
$ws.Range("H4").Value = 493.33334
$ws.Range("K4").Value = 440
$ws.Range("N4").Value = -830
$ws.Range("J4").Value = 600
$ws.Range("M4").Value = -325
$ws.Range("H132").Value = 50780
$ws.Range("N132").Value = -60900
$ws.Range("J132").Value = 50780
$ws.Range("L132").Value = 50780

# --- Sheet: CRP (85 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 40000500
$ws.Range("K6").Value = 40000500
$ws.Range("I6").Value = 40000500
$ws.Range("M6").Value = -40000387
$ws.Range("M7").Value = 21
$ws.Range("J7").Value = 380
$ws.Range("N7").Value = -606
$ws.Range("H7").Value = 284
$ws.Range("I7").Value = 92
$ws.Range("K7").Value = 92
$ws.Range("L7").Value = 380
$ws.Range("L17").Value = 60009
$ws.Range("H17").Value = 60009
$ws.Range("J17").Value = 60009
$ws.Range("N17").Value = -60357
$ws.Range("M25").ClearContents()
$ws.Range("K25").Value = 0
$ws.Range("N25").Value = -5361
$ws.Range("H25").Value = 5013
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 5013
$ws.Range("L25").Value = 5013
$ws.Range("M41").Value = -6622
$ws.Range("I41").Value = 7050
$ws.Range("J41").Value = 19000
$ws.Range("L41").Value = 19000
$ws.Range("N41").Value = -19856
$ws.Range("K41").Value = 7050
$ws.Range("H41").Value = 17008.334
$ws.Range("M50").ClearContents()
$ws.Range("L50").Value = 22055.2
$ws.Range("H50").Value = 22055.2
$ws.Range("N50").Value = -23305.2
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 22055.2
$ws.Range("K50").Value = 0
$ws.Range("L51").Value = 23732.666
$ws.Range("N51").Value = -25204.666
$ws.Range("M51").ClearContents()
$ws.Range("K51").Value = 0
$ws.Range("J51").Value = 23732.666
$ws.Range("I51").Value = 0
$ws.Range("H51").Value = 23732.666
$ws.Range("N59").Value = -34417
$ws.Range("K59").Value = 20000
$ws.Range("J59").Value = 32127
$ws.Range("I59").Value = 20000
$ws.Range("L59").Value = 32127
$ws.Range("M59").Value = -18855
$ws.Range("H59").Value = 29095.25
$ws.Range("H60").Value = 15000
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("I60").Value = 0
$ws.Range("N61").Value = -24428.666
$ws.Range("J61").Value = 23732.666
$ws.Range("H61").Value = 23732.666
$ws.Range("L61").Value = 23732.666
$ws.Range("K61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("J68").Value = 33863.332
$ws.Range("L68").Value = 33863.332
$ws.Range("H68").Value = 33863.332
$ws.Range("N68").Value = -35361.332
$ws.Range("H71").Value = 33863.332
$ws.Range("L71").Value = 101589.996
$ws.Range("N71").Value = -109077.996
$ws.Range("J71").Value = 33863.332
$ws.Range("K74").Value = 20000
$ws.Range("I74").Value = 20000
$ws.Range("H74").Value = 23333.334
$ws.Range("M74").Value = -19126
$ws.Range("H77").Value = 23333.334
$ws.Range("K77").Value = 60000
$ws.Range("M77").Value = -55632
$ws.Range("I77").Value = 20000
$ws.Range("H132").Value = 53006.85
$ws.Range("M132").Value = -5856.928400000001
$ws.Range("I132").Value = 2795.6428
$ws.Range("K132").Value = 8386.928400000001
$ws.Range("H134").Value = 99652.27
$ws.Range("I134").Value = 1335
$ws.Range("M134").Value = -1470
$ws.Range("K134").Value = 4005

# --- Sheet: CUL (33 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M11").Value = -518.8
$ws.Range("H11").Value = 208629.08
$ws.Range("K11").Value = 658.8
$ws.Range("I11").Value = 219.6
$ws.Range("I12").Value = 65
$ws.Range("J12").Value = 132.95
$ws.Range("L12").Value = 398.85
$ws.Range("H12").Value = 103.828575
$ws.Range("N12").Value = -744.8499999999999
$ws.Range("M12").Value = -22
$ws.Range("K12").Value = 195
$ws.Range("L17").Value = 3000
$ws.Range("H17").Value = 1000
$ws.Range("J17").Value = 1000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3338
$ws.Range("K17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("H122").Value = 1097.1945
$ws.Range("N122").Value = -15742.3873
$ws.Range("J122").Value = 1204.7097
$ws.Range("L122").Value = 10842.3873
$ws.Range("H132").Value = 2407.2856
$ws.Range("M132").Value = -12122.819
$ws.Range("I132").Value = 1628.091
$ws.Range("K132").Value = 14652.819
$ws.Range("K139").Value = 4374.6819
$ws.Range("L139").Value = 13500
$ws.Range("I139").Value = 1458.2273
$ws.Range("J139").Value = 4500
$ws.Range("N139").Value = -23780
$ws.Range("M139").Value = 765.3181000000004
$ws.Range("H139").Value = 1590.4783

# --- Sheet: LTW (15 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -27779922
$ws.Range("H7").Value = 22729610
$ws.Range("I7").Value = 27780034
$ws.Range("K7").Value = 27780034
$ws.Range("I126").Value = 27780034
$ws.Range("K126").Value = 83340102
$ws.Range("H126").Value = 22729610
$ws.Range("M126").Value = -83337632
$ws.Range("H132").Value = 48422.816
$ws.Range("J132").Value = 103230.4
$ws.Range("M132").Value = -5719.499899999999
$ws.Range("L132").Value = 309691.2
$ws.Range("I132").Value = 2749.8333
$ws.Range("N132").Value = -314751.2
$ws.Range("K132").Value = 8249.499899999999

# --- Sheet: WVR (21 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J81").Value = 2758.4167
$ws.Range("K81").Value = 1320
$ws.Range("L81").Value = 5516.8334
$ws.Range("H81").Value = 2141.2354
$ws.Range("N81").Value = -7638.8334
$ws.Range("I81").Value = 660
$ws.Range("M81").Value = -259
$ws.Range("L84").Value = 27584.167
$ws.Range("M84").Value = -1296
$ws.Range("N84").Value = -38192.167
$ws.Range("K84").Value = 6600
$ws.Range("I84").Value = 660
$ws.Range("J84").Value = 2758.4167
$ws.Range("H84").Value = 2141.2354
$ws.Range("H132").Value = 48087.023
$ws.Range("J132").Value = 61117.59
$ws.Range("M132").Value = -116171.12
$ws.Range("L132").Value = 183352.77
$ws.Range("I132").Value = 39567.04
$ws.Range("N132").Value = -188412.77
$ws.Range("K132").Value = 118701.12
